$p = $ppt.ActivePresentation
$EmuPerPt = 12700

# --- Slide 13: "Final Code Generator" --------------------------------
# Body placeholder (shape 4), bullet sub-items under "Compiler project"
$s13 = $p.Slides.Item(13)
$body13 = $s13.Shapes.Item(4)
$tr13 = $body13.TextFrame.TextRange
$tr13.Paragraphs(4, 1).Text = "assembler provided in the course repository performs final code generation"
$tr13.Paragraphs(5, 1).Text = "assembler also implements minor optimizations"

# --- Slide 17: "Compiler Construction Tools" --------------------------
# "x86" label inside the "Group 27" diagram -> "x86-64"
$s17 = $p.Slides.Item(17)
$group17 = $s17.Shapes.Item(5)
$x86Box = $group17.GroupItems.Item(4)
$x86Box.TextFrame.TextRange.Paragraphs(1, 1).Text = "x86-64"

# --- Slide 19: "Single-pass Versus Multi-pass Compilers" --------------
# Body placeholder (shape 4), last bullet point text tweak
$s19 = $p.Slides.Item(19)
$body19 = $s19.Shapes.Item(4)
$tr19 = $body19.TextFrame.TextRange
$tr19.Paragraphs(10, 1).Text = "requires design of intermediate languages/representations"

# Reposition the floating caption textbox (shape 5) slightly upward
$caption19 = $s19.Shapes.Item(5)
$caption19.Top = 5715000 / $EmuPerPt

# --- Slide 20: "Passes in the Compiler Project" -----------------------
$s20 = $p.Slides.Item(20)
$note20 = $s20.Shapes.Item(5)
$note20Height = $note20.Height
$tr20 = $note20.TextFrame.TextRange
$tr20.Paragraphs(2, 1).Text = "in-memory data structures called abstract syntax trees.  The only I/O to disk occurs when reading the source file and generating assembly code."
# this shape auto-fits to its text; restore the original height so only
# the wording (not the box size) changes, matching the source edit
$note20.Height = $note20Height

# --- Slide 8: "Constraint Analyzer" ------------------------------------
# Reposition the floating caption textbox (shape 5) slightly downward
$s8 = $p.Slides.Item(8)
$caption8 = $s8.Shapes.Item(5)
$caption8.Top = 4884003 / $EmuPerPt
